$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh "quantidade_atipica" (A) and "estoque_atualizado" (G) for the rows
# already present - the nightly cronjob recomputed these two figures while
# every other column (date/cliente/id_produto/produto/critico) stayed put.
$ws.Range("A2").Value = 5
$ws.Range("G2").Value = 681

$ws.Range("A3").Value = 8
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = 9
$ws.Range("G4").Value = 8432

$ws.Range("A5").Value = 1
$ws.Range("G5").Value = 196

$ws.Range("A6").Value = 4
$ws.Range("G6").Value = -108

$ws.Range("A7").Value = 0
$ws.Range("G7").Value = 215

$ws.Range("A8").Value = 2
$ws.Range("G8").Value = 2

$ws.Range("A9").Value = 6
$ws.Range("G9").Value = 489

$ws.Range("A10").Value = 7
$ws.Range("G10").Value = 1

$ws.Range("A11").Value = 10
$ws.Range("G11").Value = 150

# Append the new row for 2025-04-30 produced by this cronjob run.
# Force columns B (date text) and E (zero-padded id) to Text so Excel
# doesn't reinterpret them as a date serial / number and drop the
# formatting / leading zeros.
$ws.Range("B12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "2025-04-30"
$ws.Range("C12").Value = 250
$ws.Range("D12").Value = "RH MULTI SERVICOS ADMINISTRATIVOS S.A"
$ws.Range("E12").Value = "000041"
$ws.Range("F12").Value = "LUVAS DESCARTAVEIS C/ 100 UND"
$ws.Range("G12").Value = 1109
$ws.Range("H12").Value = $false

# Line up the new row's look with its neighbours: column A keeps the
# bold/bordered/centered label style, the rest stay plain (no explicit
# style), exactly like row 11.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
